$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.844.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.807.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.92"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4441"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3668"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07285"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8518"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.62"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.806.97"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.586"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07077"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008720"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.861.91"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.144"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.986"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.27"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.185"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.16"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08805"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7457"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.935"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.428"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.094"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.50%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5289"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.22%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.030"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1683"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5169"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +8.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.414"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.53"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.977"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.29"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.649"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06334"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9143"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.45%  "
